$d = $word.ActiveDocument

# --- Bullet 1 --------------------------------------------------------------
# Before: "Generate color patches that lie on an equi-luminant space and
#          form an ellipse (or circle) on the DKL space."
# After:  "...form an ellipse (or circle) on the DKL space (generateDKLColors
#          and generateDKLColorsOnCircle)."
$found1 = $d.Content.Find.Execute(
    "form an ellipse (or circle) on the DKL space.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "form an ellipse (or circle) on the DKL space (generateDKLColors and generateDKLColorsOnCircle).",
    2
)

# --- Bullet 2 --------------------------------------------------------------
# Before: "Generate color patches that lie on an equi-luminant plane and
#          span different saturation levels."
# After:  "...span different saturation levels going from the white point
#          towards the primaries (generateCIEpointsOnLine)."
$found2 = $d.Content.Find.Execute(
    "span different saturation levels.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "span different saturation levels going from the white point towards the primaries (generateCIEpointsOnLine).",
    2
)

Write-Output "bullet1_replaced=$found1"
Write-Output "bullet2_replaced=$found2"
